# Landscaping Data.xlsx - "Add files via upload" edit
#
# Appends 7 new observation rows (sheet rows 107-113, all dated 45802 /
# 2025-05-25) to the bottom of the single data table on Sheet1, extends the
# F-column "Temp_Diff" fill-down formula to cover the new rows, and updates
# the window/selection to reflect the newly-scrolled view (mirroring what
# Excel does automatically when a user keys in new rows at the bottom of a
# sheet and saves).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. New data rows (107-113)
# ------------------------------------------------------------------
# Column layout (row 1 headers):
#  A Date  B Plant_Type  C Plant_Size  D Low  E High  F Temp_Diff
#  G Rain  H Growth      I Pruned      J Quadrant  K Shade  L UV
#  M Humidity  N Dew_Point  O Pressure  P Wind_Gust  Q Cloud_Cover
#  R Visibility  S AQI  T Pollen

$rows = @(
    @{ Row=107; B="Flowering";    C="Large";  D=43; E=64; H=0;    I="No"; J=2; K="Bright"  },
    @{ Row=108; B="Nonflowering"; C="Medium"; D=43; E=64; H=0;    I="No"; J=3; K="Bright"  },
    @{ Row=109; B="Nonflowering"; C="Small";  D=43; E=64; H=0.1;  I="No"; J=3; K="Dark"    },
    @{ Row=110; B="Nonflowering"; C="Medium"; D=43; E=64; H=0.25; I="No"; J=3; K="Neutral" },
    @{ Row=111; B="Nonflowering"; C="Medium"; D=43; E=64; H=0.5;  I="No"; J=3; K="Neutral" },
    @{ Row=112; B="Nonflowering"; C="Large";  D=43; E=64; H=0;    I="No"; J=4; K="Dark"    },
    @{ Row=113; B="Tree";         C="Medium"; D=43; E=64; H=1.35; I="No"; J=1; K="Bright"  }
)

foreach ($r in $rows) {
    $row = $r.Row

    # Column A: date - copy the style (incl. date number format) from the
    # last existing data row so the new cells reuse the same style index
    # instead of minting a new number format, then overwrite the value.
    $ws.Cells.Item(106, 1).Copy($ws.Cells.Item($row, 1))
    $ws.Cells.Item($row, 1).Value = 45802

    $ws.Cells.Item($row, 2).Value  = $r.B
    $ws.Cells.Item($row, 3).Value  = $r.C
    $ws.Cells.Item($row, 4).Value  = $r.D
    $ws.Cells.Item($row, 5).Value  = $r.E
    # F (Temp_Diff) is filled in as a formula below, once all the source
    # rows exist.
    $ws.Cells.Item($row, 7).Value  = 0.01
    $ws.Cells.Item($row, 8).Value  = $r.H
    $ws.Cells.Item($row, 9).Value  = $r.I
    $ws.Cells.Item($row, 10).Value = $r.J
    $ws.Cells.Item($row, 11).Value = $r.K
    $ws.Cells.Item($row, 12).Value = 8
    $ws.Cells.Item($row, 13).Value = 0.5
    $ws.Cells.Item($row, 14).Value = 44
    $ws.Cells.Item($row, 15).Value = 30.22
    $ws.Cells.Item($row, 16).Value = 14
    $ws.Cells.Item($row, 17).Value = 0.73
    $ws.Cells.Item($row, 18).Value = 9.9
    $ws.Cells.Item($row, 19).Value = 34
    $ws.Cells.Item($row, 20).Value = 37
}

# ------------------------------------------------------------------
# 2. Extend the Temp_Diff fill-down formula (was F67:F106) down through
#    the newly added rows F107:F113, same as dragging the fill handle.
# ------------------------------------------------------------------
$ws.Range("F107:F113").Formula = "=ABS(D107-E107)"

# ------------------------------------------------------------------
# 3. View state: Excel persists the active window's scroll position and
#    selection on save - reflect the new bottom-of-sheet view (same
#    selection shape as before: column Q, now over the newly added rows).
# ------------------------------------------------------------------
$ws.Range("Q107:Q113").Select()
